$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the test data rows (rows 3-7), keeping header (row1) and the first
# data row (row2) intact. Only cell formatting/styles remain for rows 3-7.
$ws.Range("A3:C7").ClearContents()

# Remove the hyperlinks that were attached to the now-deleted test data
# cells, keeping only the hyperlinks on A2 and B2.
$keepAddrs = @('$A$2', '$B$2')
$changed = $true
while ($changed) {
  $changed = $false
  foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($keepAddrs -notcontains $addr) {
      $h.Delete()
      $changed = $true
      break
    }
  }
}

# Update the active selection to match the saved view state.
$ws.Range("C11").Select()
